# Updated cryptos list on Sun Mar 31 21:27:00 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.817.13'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.78%  '

# --- Row 3 ---
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.633.15'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.80%  '

# --- Row 4 ---
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '

# --- Row 5 ---
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.97'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.77%  '

# --- Row 6 ---
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '198.73'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.28%  '

# --- Row 7 ---
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.68%  '

# --- Row 9 ---
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.222'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +10.54%  '

# --- Row 10 ---
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.647'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.04%  '

# --- Row 11 ---
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.91'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.45%  '

# --- Row 12 ---
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.93%  '

# --- Row 13 ---
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.10%  '

# --- Row 14 ---
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.210.31'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.73%  '

# --- Row 15 ---
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '683.51'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +15.11%  '

# --- Row 16 ---
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.94'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.97%  '

# --- Row 17 ---
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.892.03'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.66%  '

# --- Row 18 ---
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.659.87'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.96%  '

# --- Row 19 ---
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.03'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.15%  '

# --- Row 20 ---
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.20%  '

# --- Row 21 ---
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.57%  '

# --- Row 22 ---
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.82'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.10%  '

# --- Row 23 ---
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.40'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.55%  '

# --- Row 24 ---
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '104.92'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.48%  '

# --- Row 25 ---
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.64'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.34%  '

# --- Row 26 ---
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.34%  '

# --- Row 27 ---
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.95%  '

# --- Row 28 ---
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.89'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.34%  '

# --- Row 29 ---
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.33'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.67%  '

# --- Row 30 ---
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.59'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +6.71%  '

# --- Row 31 ---
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.17'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.26%  '

# --- Row 32 ---
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.18'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.37%  '

# --- Row 33 ---
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.77%  '

# --- Row 34 ---
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.30'

# --- Row 35 ---
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0870'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.82%  '

# --- Row 36 ---
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.949.31'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +5.66%  '

# --- Row 37 ---
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.13%  '

# --- Row 38 ---
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.04'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.15%  '

# --- Row 39 ---
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.73'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.35%  '

# --- Row 40 (was Bittensor, now TheGraph) ---
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.388'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.39%  '

# --- Row 41 (was TheGraph, now Bittensor) ---
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '502.22'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.15%  '

# --- Row 42 ---
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.55'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.66%  '

# --- Row 43 ---
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.37%  '

# --- Row 44 ---
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.06'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +9.43%  '

# --- Row 45 ---
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.41%  '

# --- Row 46 ---
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.49'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +6.10%  '

# --- Row 47 ---
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.68%  '

# --- Row 48 ---
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.38%  '

# --- Row 49 ---
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.24%  '

# --- Row 50 ---
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000248'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.70%  '

# --- Row 51 ---
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.59%  '
